$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6200733351348301
$ws.Range("C3").Value = 0.6024892950496797
$ws.Range("D3").Value = 2
$ws.Range("C4").Value = 0.6287381491519984
$ws.Range("C5").Value = 0.6118046690512592
